$wb = $excel.ActiveWorkbook

# Insert a new worksheet named "RAMFlags" right after "ALU Invert Logic"
# (becomes the 2nd sheet, ahead of "Microcode" and "Opcodes").
$sheet1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "RAMFlags"

# Notes at the bottom of the sheet
$newSheet.Range("E14").Value = "R/W = clock nand (not nInput)"
$newSheet.Range("E13").Value = "OD = clock or nOutput"

# "mode" column (D) and explanatory notes column (G)
$newSheet.Range("D4").Value = "none"
$newSheet.Range("D5").Value = "none"
$newSheet.Range("D9").Value = "none"
$newSheet.Range("D8").Value = "write"
$newSheet.Range("G7").Value = "too late to read"
$newSheet.Range("D7").Value = "read/none"
$newSheet.Range("D2").Value = "undefined"
$newSheet.Range("D6").Value = "undefined"
$newSheet.Range("G4").Value = "don't want to write because address is changing"
$newSheet.Range("D3").Value = "read"
$newSheet.Range("G2").Value = "happens potentially on ROM switching"

# Column headers (row 1), entered right-to-left
$newSheet.Range("F1").Value = "R/W"
$newSheet.Range("E1").Value = "OD"
$newSheet.Range("D1").Value = "mode"
$newSheet.Range("C1").Value = "nInput"
$newSheet.Range("B1").Value = "nOutput"
$newSheet.Range("A1").Value = "clock"

# Truth-table numeric inputs/outputs
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = 0
$newSheet.Range("C2").Value = 0
$newSheet.Range("E2").Value = 1
$newSheet.Range("F2").Value = 1

$newSheet.Range("A3").Value = 0
$newSheet.Range("B3").Value = 0
$newSheet.Range("C3").Value = 1
$newSheet.Range("E3").Value = 0
$newSheet.Range("F3").Value = 1

$newSheet.Range("A4").Value = 0
$newSheet.Range("B4").Value = 1
$newSheet.Range("C4").Value = 0
$newSheet.Range("E4").Value = 1
$newSheet.Range("F4").Value = 1

$newSheet.Range("A5").Value = 0
$newSheet.Range("B5").Value = 1
$newSheet.Range("C5").Value = 1
$newSheet.Range("E5").Value = 1
$newSheet.Range("F5").Value = 1

$newSheet.Range("A6").Value = 1
$newSheet.Range("B6").Value = 0
$newSheet.Range("C6").Value = 0

$newSheet.Range("A7").Value = 1
$newSheet.Range("B7").Value = 0
$newSheet.Range("C7").Value = 1
$newSheet.Range("E7").Value = "x"
$newSheet.Range("F7").Value = 1

$newSheet.Range("A8").Value = 1
$newSheet.Range("B8").Value = 1
$newSheet.Range("C8").Value = 0
$newSheet.Range("E8").Value = "x"
$newSheet.Range("F8").Value = 0

$newSheet.Range("A9").Value = 1
$newSheet.Range("B9").Value = 1
$newSheet.Range("C9").Value = 1
$newSheet.Range("E9").Value = 1
$newSheet.Range("F9").Value = 1

$newSheet.Range("E13").Select()

$wb.Save()
